$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new draw result for 2025-11-10 as row 55
$row = 55
$ws.Cells.Item($row, 1).Value = "2025-11-10"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "251110"
$ws.Cells.Item($row, 4).Value = "1-5-5"
$ws.Cells.Item($row, 5).Value = "2025-11-10T21:37:45.601+04:00"
